# registeUser function checked successfully
# Update the sample user row on Sheet1 with a freshly-registered user's
# details (name, email, phone/"password" number) and leave the UI focused
# on Sheet1 / cell C2 (mirrors the manual check the author performed after
# running the registerUser test).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# New user details replace the old Shubham / shubham@gmail.com / shu123 row.
$ws1.Range("A2").Value = " shubham kumar"
$ws1.Range("B2").Value = "shubhamk@gmail.com"
$ws1.Range("C2").Value = 12345678

# Switch focus to Sheet1 and leave the selection on C2, like the diff shows.
[void]$ws1.Activate()
[void]$ws1.Range("C2").Select()
